# Add an "ORGANIZATION" column to the CodeSchemes sheet.
#
# Inserts a new column B (pushing ID/CLASSIFICATION/... etc. one column to
# the right), gives it the header "ORGANIZATION" and fills in the
# organization identifier for the single data row, then keeps the "yti"
# defined name (CodeSchemes!$A$1:$W$2) in sync with the now-wider range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CodeSchemes")

# Insert a new, blank column before column B. Excel copies formatting
# (e.g. column width) from the column to the left of the insertion point.
$ws.Columns("B").Insert()
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

# New header + value for the inserted column.
$ws.Range("B1").Value = "ORGANIZATION"
$ws.Range("B2").Value = "74a41211-8c99-4835-a519-7a61612b1098"

# The "yti" defined name pointed at CodeSchemes!$A$1:$W$2; now that the
# table grew by one column it should cover through column X.
$name = $wb.Names.Item("yti")
$name.RefersTo = "=CodeSchemes!`$A`$1:`$X`$2"
